$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.744.61"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.630.18"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.43"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0633"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.53"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.853.02"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "1.625.44"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.554"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "0.0₃0761"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.84"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "25.733.34"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.53"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.92"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.22"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0493"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "1.139.64"
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.72"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.54"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.805"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "1.764.21"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.14"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0509"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.45"
$ws.Range("E50").Value = "  +5.57%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").Value = "  +1.14%  "
